$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 53, shifting existing rows 53-187 down to 54-188
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new data
$ws.Cells.Item(53, 1).Value = 8
$ws.Cells.Item(53, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(53, 3).Value = "Coquimbo"
$ws.Cells.Item(53, 4).Value = 44690
$ws.Cells.Item(53, 5).Value = 4
$ws.Cells.Item(53, 6).Value = 100112037
$ws.Cells.Item(53, 7).Value = "Cebollín"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 800
$ws.Cells.Item(53, 11).Value = 1100
$ws.Cells.Item(53, 12).Value = 1200
$ws.Cells.Item(53, 13).Value = 1150
$ws.Cells.Item(53, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(53, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(53, 16).Value = 192
$ws.Cells.Item(53, 17).Value = 6
$ws.Cells.Item(53, 18).Value = "Hortaliza"
